## Aragon hospital COVID dataset — append the 2020-06-03 (serial 43985) block
## of 20 hospital rows, mirroring the previous day's (43984) block of rows
## 1193:1212 both in layout/formatting and in the static columns
## (hospital/municipio/provincia/codigo_ine/observaciones), while updating the
## date and the two occupancy counters (camas_uci_ocupadas, camas_ocupadas_total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the prior day's 20-row block (values first, then formats) into
#    the new block immediately below the existing data. Column A (the date
#    serial) is copied separately via direct cell assignment below, since
#    pasting values of that date-formatted source column makes the host
#    register an extra (unused) number format / cell style that the target
#    workbook never had.
$ws.Range("B1193:H1212").Copy()
$ws.Range("B1213:H1232").PasteSpecial(-4163)   # xlPasteValues

# 2) New date for the whole block: 2020-06-03
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item(1213 + $i, 1).Value = 43985
}

$ws.Range("A1193:H1212").Copy()
$ws.Range("A1213:H1232").PasteSpecial(-4122)   # xlPasteFormats

# 3) Updated occupancy counts (camas_uci_ocupadas = C, camas_ocupadas_total = D)
$ws.Cells.Item(1213, 3).Value = 6
$ws.Cells.Item(1213, 4).Value = 4

$ws.Cells.Item(1214, 3).Value = 31

$ws.Cells.Item(1215, 3).Value = 6

$ws.Cells.Item(1216, 3).Value = 1

$ws.Cells.Item(1217, 3).Value = 2

$ws.Cells.Item(1218, 3).Value = 4
$ws.Range("D1218").ClearContents()

$ws.Cells.Item(1219, 3).Value = 2

$ws.Cells.Item(1220, 3).Value = 11
$ws.Cells.Item(1220, 4).Value = 1

$ws.Cells.Item(1221, 3).Value = 7
$ws.Cells.Item(1221, 4).Value = 1

# row 1222 (Hospital Sagrado Corazón): C/D stay blank

$ws.Cells.Item(1223, 3).Value = 3

$ws.Cells.Item(1224, 3).Value = 6

# row 1225 (Hospital Ejea - Cinco Villas): C/D stay blank

# row 1226 (MAZ): C/D stay blank

$ws.Cells.Item(1227, 3).Value = 1

$ws.Range("C1228").ClearContents()

$ws.Cells.Item(1229, 3).Value = 1

# row 1230 (Hospital San Juan de Dios de Zaragoza): C/D stay blank
# row 1231 (Clínica Viamed Santiago): C/D stay blank
# row 1232 (Clínica El Pilar): C/D stay blank
